# Refresh the cryptos list (Price / Volume(1h) columns) with updated figures.
# Price-column values are leading-apostrophe-prefixed so Excel's COM Value
# setter stores them as text (matching the workbook's original inline-string
# cells) instead of silently coercing look-alike numbers (e.g. "9.770",
# "1.001") into doubles and dropping formatting such as trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.481.87"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "'1.573.04"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "'292.04"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'0.3724"
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("D8").Value = "'49.99"
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").Value = "'0.3401"
$ws.Range("D10").Value = "'1.146"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Value = "'0.07559"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "'21.33"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").Value = "'6.049"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").Value = "'6.974"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").Value = "'1.573.01"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "'90.83"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "'0.06767"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "'6.302"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").Value = "'16.38"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").Value = "'12.21"
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("D24").Value = "'22.486.01"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").Value = "'2.369"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").Value = "'20.04"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").Value = "'149.32"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").Value = "'5.053"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").Value = "'125.32"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").Value = "'1.748.98"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("E32").Value = "  +8.88%  "
$ws.Range("D33").Value = "'6.258"
$ws.Range("E33").Value = "  +2.62%  "
$ws.Range("D34").Value = "'2.008"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").Value = "'9.770"
$ws.Range("E35").Value = "  -3.43%  "
$ws.Range("D36").Value = "'0.08362"
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("D37").Value = "'0.02490"
$ws.Range("E37").Value = "  -2.06%  "
$ws.Range("D38").Value = "'0.2307"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("E39").Value = "  -3.70%  "
$ws.Range("D40").Value = "'0.06523"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("D41").Value = "'5.466"
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "'0.6251"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "'14.02"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "'3.812"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").Value = "'0.5876"
$ws.Range("D48").Value = "'130.73"
$ws.Range("E48").Value = "  +5.03%  "
$ws.Range("D49").Value = "'2.076"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").Value = "'1.215"
$ws.Range("E50").Value = "  -5.32%  "
$ws.Range("D51").Value = "'0.07337"
$ws.Range("E51").Value = "  +0.26%  "
